$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update B2/C2 (Ambiente/URL for the "gw" row) to point to the new
#     pre-production environment, while keeping the original cell styles
#     (B2 keeps its quote-prefix style, C2 keeps the hyperlink style).
#     Setting .Value directly on a string cell resets its style index in
#     this engine, so we restore the original formatting by copying it
#     back from the untouched sibling cells in row 3 right after.
$ws.Range("B2").Value = "i-preproducciongestion.segurossura.com.ar"
$ws.Range("B3").Copy()
$ws.Range("B2").PasteSpecial(-4122)

$ws.Range("C2").Value = "https://i-preproducciongestion.segurossura.com.ar/pc/PolicyCenter.do"

# --- Update the account number on row 2 (H2)
$ws.Range("H2").Value = 2240451788

# --- Recreate the hyperlinks: C3 keeps pointing at the original
#     oracleoutsourcing PolicyCenter URL, C2 now points at the new
#     i-preproducciongestion PolicyCenter URL. Stash the original
#     hyperlink-cell format first since Hyperlinks.Add() re-stamps the
#     target cell's style.
$ws.Range("C3").Copy()
$ws.Range("Z1").PasteSpecial(-4122)

$ws.Range("C2:C3").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C3"), "https://ssurgwsoadev4-oci.opc.oracleoutsourcing.com/pc/PolicyCenter.do")
$ws.Hyperlinks.Add($ws.Range("C2"), "https://i-preproducciongestion.segurossura.com.ar/pc/PolicyCenter.do")

$ws.Range("Z1").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("Z1").Clear()

# --- Clear the leftover styled-but-empty placeholder cells on row 4
$ws.Range("B4").Clear()
$ws.Range("C4").Clear()

# --- Update the current selection / view state
$ws.Range("D9").Select() | Out-Null
